$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 34

# Columns A-D hold text that looks numeric/date/time-like (dates, times,
# weekday names, and a plain "23"). Force them to be treated as literal
# text (inline/shared string) rather than being auto-converted by Excel
# into date/time serial numbers, then strip the temporary text format so
# the cells end up unstyled, same as the rest of the sheet.
$textRange = $ws.Range("A${row}:D${row}")
$textRange.NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "2023-06-08"
$ws.Cells.Item($row, 2).Value = "22:27:44"
$ws.Cells.Item($row, 3).Value = "Thursday"
$ws.Cells.Item($row, 4).Value = "23"

$textRange.ClearFormats()

$ws.Cells.Item($row, 5).Value = 120109
$ws.Cells.Item($row, 6).Value = 134373
$ws.Cells.Item($row, 7).Value = 160557
$ws.Cells.Item($row, 8).Value = 131487
$ws.Cells.Item($row, 9).Value = 175611
$ws.Cells.Item($row, 10).Value = 113316
$ws.Cells.Item($row, 11).Value = 201350
$ws.Cells.Item($row, 12).Value = 221446
$ws.Cells.Item($row, 13).Value = 173060
$ws.Cells.Item($row, 14).Value = 120177
$ws.Cells.Item($row, 15).Value = 38670
$ws.Cells.Item($row, 16).Value = 34427
$ws.Cells.Item($row, 17).Value = 50891
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 36975
$ws.Cells.Item($row, 20).Value = -1
